$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 updates
$ws.Range("I3").Value = "Segunda"
$ws.Range("J3").Value = 60
$ws.Range("K3").Value = 15000
$ws.Range("L3").Value = 15000
$ws.Range("M3").Value = 15000
$ws.Range("P3").Value = 833

# Row 4 updates
$ws.Range("D4").Value = 44637
$ws.Range("I4").Value = "Primera"
$ws.Range("J4").Value = 100
$ws.Range("L4").Value = 16000
$ws.Range("M4").Value = 15500
$ws.Range("P4").Value = 861

# Row 5 updates
$ws.Range("D5").Value = 44656

# Row 6 updates
$ws.Range("D6").Value = 44658
$ws.Range("J6").Value = 80
